$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New card/effect under the "CARDS" section (row 2, col P)
$ws.Range("P2").Value = "Excite"

# New "Other cards" block (rows 10-12)
$ws.Range("J10").Value = "Other cards"

$ws.Range("J11").Value = "Force"
$ws.Range("K11").Value = "Force a bid"

$ws.Range("J12").Value = "Skip"
$ws.Range("K12").Value = "Skip a bid"
$ws.Range("L12").Value = "Reduce fatigure"

# Leave the selection where the author left off editing
$ws.Range("L12").Select()
